# Updates the '想去人数' (want-to-go count) column F across sheets,
# matching gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1321
$ws.Range("F3").Value = 77
$ws.Range("F4").Value = 74
$ws.Range("F5").Value = 165
$ws.Range("F6").Value = 386
$ws.Range("F7").Value = 172
$ws.Range("F8").Value = 121
$ws.Range("F9").Value = 987
$ws.Range("F10").Value = 318
$ws.Range("F11").Value = 174
$ws.Range("F14").Value = 358
$ws.Range("F15").Value = 343
$ws.Range("F16").Value = 750
$ws.Range("F17").Value = 127
$ws.Range("F18").Value = 702
$ws.Range("F19").Value = 248
$ws.Range("F20").Value = 63
$ws.Range("F21").Value = 965
$ws.Range("F22").Value = 429
$ws.Range("F23").Value = 239
$ws.Range("F24").Value = 76
$ws.Range("F25").Value = 351
$ws.Range("F28").Value = 451

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 32
$ws.Range("F11").Value = 144

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1321
$ws.Range("F4").Value = 77
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 165
$ws.Range("F8").Value = 386
$ws.Range("F9").Value = 172
$ws.Range("F10").Value = 121
$ws.Range("F11").Value = 987
$ws.Range("F12").Value = 318
$ws.Range("F13").Value = 174
$ws.Range("F18").Value = 32
$ws.Range("F19").Value = 358
$ws.Range("F22").Value = 343
$ws.Range("F23").Value = 750
$ws.Range("F24").Value = 127
$ws.Range("F25").Value = 702
$ws.Range("F26").Value = 248
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 965
$ws.Range("F29").Value = 429
$ws.Range("F32").Value = 239
$ws.Range("F33").Value = 76
$ws.Range("F34").Value = 351
$ws.Range("F36").Value = 144
$ws.Range("F40").Value = 451
